$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.272.52"
$ws.Range("E2").Value = "  +7.12%  "
$ws.Range("D3").Value = "3.883.65"
$ws.Range("E3").Value = "  +15.41%  "
$ws.Range("E4").Value = "  -1.37%  "
$ws.Range("D5").Value = "'425.81"
$ws.Range("E5").Value = "  +10.77%  "
$ws.Range("D6").Value = "'131.35"
$ws.Range("E6").Value = "  +9.03%  "
$ws.Range("D7").Value = "3.875.94"
$ws.Range("E7").Value = "  +9.46%  "
$ws.Range("D8").Value = "'0.615"
$ws.Range("E8").Value = "  +7.51%  "
$ws.Range("D9").Value = "'0.997"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  +11.65%  "
$ws.Range("E11").Value = "  +13.59%  "
$ws.Range("D12").Value = "'0.0000347"
$ws.Range("E12").Value = "  +19.52%  "
$ws.Range("D13").Value = "'41.04"
$ws.Range("E13").Value = "  +8.42%  "
$ws.Range("D14").Value = "4.484.76"
$ws.Range("E14").Value = "  +14.08%  "
$ws.Range("E15").Value = "  +14.42%  "
$ws.Range("D16").Value = "'15.96"
$ws.Range("E16").Value = "  +31.01%  "
$ws.Range("D17").Value = "3.883.17"
$ws.Range("E17").Value = "  +12.22%  "
$ws.Range("D18").Value = "'0.138"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("E19").Value = "  +10.17%  "
$ws.Range("D20").Value = "67.204.23"
$ws.Range("E20").Value = "  +6.71%  "
$ws.Range("E21").Value = "  +8.63%  "
$ws.Range("D22").Value = "'414.56"
$ws.Range("E22").Value = "  +10.56%  "
$ws.Range("D23").Value = "'14.91"
$ws.Range("E23").Value = "  +11.84%  "
$ws.Range("D24").Value = "'84.45"
$ws.Range("E24").Value = "  +7.54%  "
$ws.Range("E25").Value = "  +10.68%  "
$ws.Range("D26").Value = "'37.68"
$ws.Range("E26").Value = "  +16.17%  "
$ws.Range("D27").Value = "'9.96"
$ws.Range("E27").Value = "  +16.79%  "
$ws.Range("D28").Value = "'3.25"
$ws.Range("E28").Value = "  +11.52%  "
$ws.Range("D29").Value = "'5.31"
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("D30").Value = "'8.96"
$ws.Range("E30").Value = "  +40.17%  "
$ws.Range("D31").Value = "'733.82"
$ws.Range("E31").Value = "  +13.81%  "
$ws.Range("D32").Value = "'13.38"
$ws.Range("E32").Value = "  +15.29%  "
$ws.Range("E33").Value = "  +15.03%  "
$ws.Range("E34").Value = "  +7.71%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'39.05"
$ws.Range("E36").Value = "  +7.81%  "
$ws.Range("E37").Value = "  +4.85%  "
$ws.Range("D38").Value = "'55.75"
$ws.Range("E38").Value = "  +4.11%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0768"
$ws.Range("E39").Value = "  +31.32%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").Value = "'5.31"
$ws.Range("E40").Value = "  +36.86%  "
$ws.Range("E41").Value = "  +8.96%  "
$ws.Range("D42").Value = "'2.88"
$ws.Range("E42").Value = "  +8.92%  "
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("E44").Value = "  +12.75%  "
$ws.Range("E45").Value = "  +4.64%  "
$ws.Range("D46").Value = "'3.16"
$ws.Range("E46").Value = "  +6.10%  "
$ws.Range("E47").Value = "  +18.88%  "
$ws.Range("E48").Value = "  +8.23%  "
$ws.Range("D49").Value = "'141.68"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("E50").Value = "  +9.80%  "
$ws.Range("E51").Value = "  +6.53%"
